# "Adjustment CUM/KWh - PJC"
#
# The sheet had a duplicate-ish pair of header labels: "KWH" (used as the
# units header in H11) and "Total KWh" / "KWh" appearing elsewhere. The
# author retyped the units header in H11 from "KWH" to "KWh", which (because
# of how the shared-string table gets rebuilt on save) reshuffles the
# shared-string indices used by the other label cells even though their
# visible text doesn't change. The author also re-selected cell A9 (the
# "Total KWh" label) which picked up the same cell style already used by the
# other left-column labels (A2:A8), and left the active selection on C15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Units header in the table: "KWH" -> "KWh"
$ws.Range("H11").Value = "KWh"

# Normalize the "Total KWh" label (A9) onto the same style used by the other
# left-hand labels (A2:A8), matching what Excel applied when the cell was
# touched.
$ws.Range("A9").Style = "Normal 2"
$ws.Range("A9").Font.Bold = $true

# Leave the cursor where the author left it.
[void]$ws.Range("C15").Select()
